# issue #5: stock data from json to db
# Add "category", "source_file" and "index" columns to the 股票 (stock) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = 18

# Insert a new column I ("category") between "property_category" (H) and
# "date" (old I, now shifts to J). This also shifts legislator_name/
# legislator_id from J/K to K/L. Inserting (rather than just writing into
# blank cells) copies the neighbouring cell formatting (s="1"/s="2").
$ws.Columns.Item(9).Insert()

# Insert two more new (blank) columns at the end for "source_file" (M) and
# "index" (N), again via Insert() so they pick up matching cell styles.
$ws.Columns.Item(13).Insert()
$ws.Columns.Item(14).Insert()

# Header row
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Data rows: category = "normal", source_file = "tmp93a21", index = same
# value as column A on that row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
    $ws.Cells.Item($r, 13).Value = "tmp93a21"
    $idx = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 14).Value = $idx
}
